# Apply updated classification-report metrics (BERT embedding extraction results).
# Row labels in column A are re-ordered/re-mapped and the precision/recall/f1-score/support
# values (columns B-E) are refreshed to reflect the new run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: '50' -> '50'
$ws.Range("B2").Value = 0.8
$ws.Range("D2").Value = 0.8000000000000002

# Row 3: '14' -> '14'
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 0.7499999999999999

# Row 4: '8' -> '13'
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "13"
$ws.Range("A4").NumberFormat = "General"
$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 0.6666666666666666
$ws.Range("D4").Value = 0.7272727272727272

# Row 6: '0' -> '0'
$ws.Range("B6").Value = 0.8125
$ws.Range("C6").Value = 0.9285714285714286
$ws.Range("D6").Value = 0.8666666666666666

# Row 7: '53' -> '8'
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "8"
$ws.Range("A7").NumberFormat = "General"
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 1

# Row 8: '42' -> '42'
$ws.Range("C8").Value = 0.6666666666666666
$ws.Range("D8").Value = 0.8

# Row 9: '19' -> '53'
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "53"
$ws.Range("A9").NumberFormat = "General"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.8
$ws.Range("D9").Value = 0.888888888888889

# Row 10: '46' -> '19'
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "19"
$ws.Range("A10").NumberFormat = "General"

# Row 11: '44' -> '44'
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 1

# Row 12: '39' -> '33'
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "33"
$ws.Range("A12").NumberFormat = "General"
$ws.Range("B12").Value = 0.7142857142857143
$ws.Range("C12").Value = 0.8333333333333334
$ws.Range("D12").Value = 0.7692307692307692
$ws.Range("E12").Value = 6

# Row 13: '2' -> '46'
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "46"
$ws.Range("A13").NumberFormat = "General"
$ws.Range("B13").Value = 0.5
$ws.Range("C13").Value = 0.5
$ws.Range("D13").Value = 0.5

# Row 14: '15' -> '39'
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "39"
$ws.Range("A14").NumberFormat = "General"
$ws.Range("C14").Value = 0.3333333333333333
$ws.Range("D14").Value = 0.5

# Row 15: '5' -> '21'
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "21"
$ws.Range("A15").NumberFormat = "General"

# Row 16: '11' -> '2'
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2"
$ws.Range("A16").NumberFormat = "General"

# Row 17: '16' -> '15'
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "15"
$ws.Range("A17").NumberFormat = "General"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1

# Row 18: '28' -> '5'
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "5"
$ws.Range("A18").NumberFormat = "General"
$ws.Range("B18").Value = 0.6666666666666666
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0.8

# Row 19: '27' -> '11'
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "11"
$ws.Range("A19").NumberFormat = "General"
$ws.Range("B19").Value = 0.75
$ws.Range("D19").Value = 0.8571428571428571

# Row 20: '32' -> '32'
$ws.Range("B20").Value = 0.6666666666666666
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.8

# Row 21: '33' -> '16'
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "16"
$ws.Range("A21").NumberFormat = "General"
$ws.Range("B21").Value = 0.6666666666666666
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0.8

# Row 22: '21' -> '28'
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "28"
$ws.Range("A22").NumberFormat = "General"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0.6666666666666666
$ws.Range("D22").Value = 0.8

# Row 23: '18' -> '27'
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "27"
$ws.Range("A23").NumberFormat = "General"
$ws.Range("C23").Value = 0.75
$ws.Range("D23").Value = 0.8571428571428571

# Row 24: '13' -> '18'
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "18"
$ws.Range("A24").NumberFormat = "General"

# Row 25: 'accuracy' -> 'accuracy'
$ws.Range("B25").Value = 0.8390804597701149
$ws.Range("C25").Value = 0.8390804597701149
$ws.Range("D25").Value = 0.8390804597701149
$ws.Range("E25").Value = 0.8390804597701149

# Row 26: 'macro avg' -> 'macro avg'
$ws.Range("B26").Value = 0.8642080745341615
$ws.Range("C26").Value = 0.849792960662526
$ws.Range("D26").Value = 0.8361143687230645
$ws.Range("E26").Value = 87

# Row 27: 'weighted avg' -> 'weighted avg'
$ws.Range("B27").Value = 0.8644909688013137
$ws.Range("C27").Value = 0.8390804597701149
$ws.Range("D27").Value = 0.8339040908006425
$ws.Range("E27").Value = 87
